# Update the "Förändrad" (Changed) date column (C) from 45192 to 45202
# for every data row (rows 2 through 204) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Range("A1").End(4).Row  # xlDown = 4, expand to bottom of column A

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45192) {
        $cell.Value = 45202
    }
}
